$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Relais 1pool / Conrad): update hyperlink display URL to the 24V version
# (edited first so it lands earlier in the shared-strings table, matching the source edit order)
$ws.Range("E7").Value = "https://www.conrad.be/p/finder-403170240000-printrelais-24-vdc-12-a-1x-wisselcontact-1-stuks-1560602"

# Row 6 (Relais 2-pool / Conrad): update hyperlink display URL to the 24V version
$ws.Range("E6").Value = "https://www.conrad.be/p/finder-405290240000-printrelais-24-vdc-8-a-2x-wisselcontact-1-stuks-502882"

# New row 12: Fly back diode, supplied by Mouser (no price/qty/URL yet)
$ws.Range("A12").Value = "Fly back diode"
$ws.Range("D12").Value = "Mouser"

# Restore the selection to match the saved view state
$ws.Range("B17").Select()
